$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRows = @(
    @("Laois Hire O'Moore Park", "Laois"),
    @("Heartland Credit Union Páirc Seán MacDiarmada", "Leitrim"),
    @("Fitzgerald Stadium, Killarney", "Kerry"),
    @("SuperValu Páirc Uí Chaoimh", "Cork"),
    @("Brewster Park", "Fermanagh"),
    @("TEG Cusack Park", "Westmeath"),
    @("Roger Casements, Portglenone", "Antrim"),
    @("King & Moffatt Dr. Hyde Park", "Roscommon")
)

$startRow = 113
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($row, 2).Value = $newRows[$i][1]
}
